{"js": "// Remove the \"IT Support Intern\" text that was filled in under the\n// \"Department Assigned:\" table cell, leaving the paragraph (and its\n// formatting) empty, and remove the leftover \"_GoBack\" bookmark that\n// Word inserts at the last edit position.\n\n// 1) Clear the \"IT Support Intern\" run from the table cell.\nconst body = context.document.body;\nconst results = body.search(\"IT Support Intern\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replacing with an empty string removes the run's text while\n  // preserving the (now empty) paragraph and its paragraph properties.\n  results.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the \"_GoBack\" bookmark left over from the last edit location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Remove the \"IT Support Intern\" text that was filled in under the\n# \"Department Assigned:\" table cell, leaving the paragraph (and its\n# formatting) empty, and remove the leftover \"_GoBack\" bookmark that\n# Word inserts at the last edit position.\n\n$d = $word.ActiveDocument\n\n# 1) Clear the \"IT Support Intern\" run from the table cell.\n$range = $d.Content\nif ($range.Find.Execute(\"IT Support Intern\")) {\n    $range.Text = \"\"\n}\n\n# 2) Remove the \"_GoBack\" bookmark left over from the last edit location.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
